# Append incident-log rows 16-30 (2024-05-13 data) to Sheet1, extending the
# populated range from A1:G15 to A1:G30, matching the source data exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("A16").Value = "'2024-05-13"
$ws.Range("B16").Value = '11:37:38'
$ws.Range("C16").Value = '-'
$ws.Range("D16").Value = '-'
$ws.Range("E16").Value = '-'
$ws.Range("F16").Value = 'Robot no coloca bien filter en palet'
$ws.Range("G16").Value = '-'

# Row 17
$ws.Range("A17").Value = "'2024-05-13"
$ws.Range("B17").Value = '11:45:06'
$ws.Range("C17").Value = '-'
$ws.Range("D17").Value = 'Tornillo atascado en tolva'
$ws.Range("E17").Value = '-'
$ws.Range("F17").Value = '-'
$ws.Range("G17").Value = '-'

# Row 18
$ws.Range("A18").Value = "'2024-05-13"
$ws.Range("B18").Value = '11:47:46'
$ws.Range("C18").Value = 'No atornilla tapa'
$ws.Range("D18").Value = '-'
$ws.Range("E18").Value = '-'
$ws.Range("F18").Value = '-'
$ws.Range("G18").Value = '-'

# Row 19
$ws.Range("A19").Value = "'2024-05-13"
$ws.Range("B19").Value = '11:52:36'
$ws.Range("C19").Value = '-'
$ws.Range("D19").Value = '-'
$ws.Range("E19").Value = '-'
$ws.Range("F19").Value = 'Traza'
$ws.Range("G19").Value = '-'

# Row 20
$ws.Range("A20").Value = "'2024-05-13"
$ws.Range("B20").Value = '11:52:39'
$ws.Range("C20").Value = '-'
$ws.Range("D20").Value = '-'
$ws.Range("E20").Value = '-'
$ws.Range("F20").Value = 'Fallo visión core'
$ws.Range("G20").Value = '-'

# Row 21
$ws.Range("A21").Value = "'2024-05-13"
$ws.Range("B21").Value = '11:52:45'
$ws.Range("C21").Value = 'Ascensor no sube'
$ws.Range("D21").Value = '-'
$ws.Range("E21").Value = '-'
$ws.Range("F21").Value = '-'
$ws.Range("G21").Value = '-'

# Row 22
$ws.Range("A22").Value = "'2024-05-13"
$ws.Range("B22").Value = '11:52:50'
$ws.Range("C22").Value = 'Secuencia atornillador'
$ws.Range("D22").Value = '-'
$ws.Range("E22").Value = '-'
$ws.Range("F22").Value = '-'
$ws.Range("G22").Value = '-'

# Row 23
$ws.Range("A23").Value = "'2024-05-13"
$ws.Range("B23").Value = '11:52:55'
$ws.Range("C23").Value = '-'
$ws.Range("D23").Value = '-'
$ws.Range("E23").Value = '-'
$ws.Range("F23").Value = '-'
$ws.Range("G23").Value = 'Colisión placas'

# Row 24
$ws.Range("A24").Value = "'2024-05-13"
$ws.Range("B24").Value = '11:53:03'
$ws.Range("C24").Value = 'Ascensor no sube'
$ws.Range("D24").Value = '-'
$ws.Range("E24").Value = '-'
$ws.Range("F24").Value = '-'
$ws.Range("G24").Value = '-'

# Row 25
$ws.Range("A25").Value = "'2024-05-13"
$ws.Range("B25").Value = '11:53:06'
$ws.Range("C25").Value = 'Fallo en paletizador'
$ws.Range("D25").Value = '-'
$ws.Range("E25").Value = '-'
$ws.Range("F25").Value = '-'
$ws.Range("G25").Value = '-'

# Row 26
$ws.Range("A26").Value = "'2024-05-13"
$ws.Range("B26").Value = '11:53:35'
$ws.Range("C26").Value = '-'
$ws.Range("D26").Value = 'Tornillo atascado en tolva'
$ws.Range("E26").Value = '-'
$ws.Range("F26").Value = '-'
$ws.Range("G26").Value = '-'

# Row 27
$ws.Range("A27").Value = "'2024-05-13"
$ws.Range("B27").Value = '11:53:39'
$ws.Range("C27").Value = '-'
$ws.Range("D27").Value = 'Cámara no detecta Top cover'
$ws.Range("E27").Value = '-'
$ws.Range("F27").Value = '-'
$ws.Range("G27").Value = '-'

# Row 28
$ws.Range("A28").Value = "'2024-05-13"
$ws.Range("B28").Value = '11:55:00'
$ws.Range("C28").Value = 'No pone tornillo'
$ws.Range("D28").Value = '-'
$ws.Range("E28").Value = '-'
$ws.Range("F28").Value = '-'
$ws.Range("G28").Value = '-'

# Row 29
$ws.Range("A29").Value = "'2024-05-13"
$ws.Range("B29").Value = '11:55:04'
$ws.Range("C29").Value = '-'
$ws.Range("D29").Value = '-'
$ws.Range("E29").Value = 'Tornillo atascado'
$ws.Range("F29").Value = '-'
$ws.Range("G29").Value = '-'

# Row 30
$ws.Range("A30").Value = "'2024-05-13"
$ws.Range("B30").Value = '11:55:07'
$ws.Range("C30").Value = '-'
$ws.Range("D30").Value = '-'
$ws.Range("E30").Value = '-'
$ws.Range("F30").Value = '-'
$ws.Range("G30").Value = 'Colisión placas'

# Strip the temporary "quote prefix" formatting picked up by column A above,
# restoring the default (General) cell style/number format.
$ws.Range("A16:A30").ClearFormats()
